# Weekly update: a new price record for Albahaca (Feria Lagunitas de Puerto
# Montt) is prepended at row 22, pushing the existing rows 22-118 down to
# 23-119 (dimension grows from A1:R118 to A1:R119).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 22..118 down by one, leaving a fresh blank row 22.
$ws.Rows.Item(22).Insert()

# Populate the new row with the latest weekly record.
$ws.Cells.Item(22, 1).Value = 4
$ws.Cells.Item(22, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(22, 3).Value = "Los Lagos"
$ws.Cells.Item(22, 4).Value = "2022-05-06"
$ws.Cells.Item(22, 5).Value = 10
$ws.Cells.Item(22, 6).Value = 100112052
$ws.Cells.Item(22, 7).Value = "Albahaca"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 11).Value = 8000
$ws.Cells.Item(22, 12).Value = 8000
$ws.Cells.Item(22, 13).Value = 8000
$ws.Cells.Item(22, 14).Value = "`$/docena de matas"
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(22, 16).Value = 1333
$ws.Cells.Item(22, 17).Value = 6
$ws.Cells.Item(22, 18).Value = "Hortaliza"
